$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 0.5 hours on Thursday (G) for "Sponsor Work" row 13
$ws.Range("G13").Value = 0.5

# Update Daily Total for Sponsor Work row (I13): 1 -> 1.5
$ws.Range("I13").Value = 1.5

# Add 0.5 hours on Thursday (G) for "Daily Total" row 14
$ws.Range("G14").Value = 0.5

# Update Weekly Total for Thursday/overall (I14): 6 -> 6.5
$ws.Range("I14").Value = 6.5

# Update the active selection to K6
$ws.Range("K6").Select()
